$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header changes
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Row 2 data changes
$ws.Range("B2").Value = 6.5925963874547371
$ws.Range("C2").ClearContents()
$ws.Range("D2").Value = 4.0538500496436942
$ws.Range("E2").ClearContents()

# Row 3 data changes
$ws.Range("B3").Value = 5.4026204843154222
$ws.Range("C3").Value = 6.3751365426387139
$ws.Range("D3").Value = 3.5345482465571889
$ws.Range("E3").Value = 8.099961900979336

# Update selection to match new sqref
$ws.Range("B1:E3").Select()
